$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29.953821277005368
$ws.Range("B3").Value = 2.8282432808873317
$ws.Range("B4").Value = 42.331879980233133
$ws.Range("B5").Value = 24.850481699987352
$ws.Range("B6").Value = 5.2102380718459118
$ws.Range("B7").Value = 41.76100252053034
$ws.Range("B8").Value = 29.441554371981745
$ws.Range("B9").Value = 30.693010113113001
